# edit.ps1 - apply the "Add files via upload" revision to curtisYungenResume.docx
#
# Summary of changes:
#   1. Programming Skills line: "CSS3" -> "CSS, Sass"; insert "TypeScript, " after
#      "JavaScript, "; insert "REST " before "APIs".
#   2. Remove the entire "Conway's Game of Life" personal-project entry
#      (heading/link paragraph + its two bullet paragraphs).
#   3. Move the hidden "_GoBack" bookmark from the end of the "...for version
#      control." bullet to the now-empty spacer paragraph left behind by the
#      removed Conway's Game of Life block.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Programming Skills line
# ---------------------------------------------------------------------------
# Locate the paragraph that contains the skills list so this keeps working
# even if paragraph numbering shifts a little.
$skillsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "HTML5, CSS3*") {
        $skillsPara = $cand
        break
    }
}

if ($skillsPara -eq $null) {
    throw "Could not locate the Programming Skills paragraph"
}

$skillsStart = $skillsPara.Range.Start
$skillsEnd = $skillsPara.Range.End - 1   # exclude the paragraph mark
$skillsRange = $d.Range($skillsStart, $skillsEnd)

$skillsXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>HTML</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve">5, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>CSS</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve"> Sass,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>JavaScript</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>TypeScript</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>j</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>Query</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve"> SQL, Bootstrap,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>Express.js, ReactJS, Node.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>js</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>Redux</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>Heroku</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve">REST </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Georgia" w:hAnsi="Cambria" w:cs="Calibri"/><w:color w:val="auto"/></w:rPr><w:t>APIs</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

$skillsRange.InsertXML($skillsXml)

# ---------------------------------------------------------------------------
# 2. Remove the "Conway's Game of Life" project block
# ---------------------------------------------------------------------------
$conwayStart = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Conway*Game of Life*") {
        $conwayStart = $i
        break
    }
}

if ($conwayStart -eq $null) {
    throw "Could not locate the Conway's Game of Life paragraph"
}

# The project entry spans 3 paragraphs: the title/link line and two bullets.
$pFirst = $d.Paragraphs.Item($conwayStart)
$pLast = $d.Paragraphs.Item($conwayStart + 2)
$conwayRange = $d.Range($pFirst.Range.Start, $pLast.Range.End)
$conwayRange.Delete()

# ---------------------------------------------------------------------------
# 3. Move the hidden "_GoBack" bookmark
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# The empty spacer paragraph left behind by the removed block is now at the
# same index the Conway's Game of Life title used to occupy.
$spacerPara = $d.Paragraphs.Item($conwayStart)
$d.Bookmarks.Add("_GoBack", $spacerPara.Range)

Write-Host "Edit complete."
